# Apply the edits described in the commit "last version in laptop"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The picture-path label in A5 had a casing fix: "pic/EndExp.PNG" -> "pic/EndExp.png"
# (re-saving this changed string also naturally reorders/dedupes the shared-strings
# table the same way Excel does, moving the edited entry to the end)
$ws.Range("A5").Value = "pic/EndExp.png"

# Column A was given an explicit width of 16 characters
$ws.Columns.Item(1).ColumnWidth = 15.25

# The active/selected cell on the sheet moved from F6 to H6
$ws.Range("H6").Select()

# The workbook window was resized/repositioned on screen
$win = $excel.ActiveWindow
$win.Left = 2340
$win.Top = 2340
$win.Width = 21600
$win.Height = 11295
